$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete rows 13-17 (data now only goes through row 12)
$ws.Rows("13:17").Delete()

# 2) Add header values for the new columns I, J, K (copy style from H1 so
#    the new header cells share style index 1, same as B1:H1)
$ws.Range("H1").Copy()
$ws.Range("I1:K1").PasteSpecial(-4122)
$ws.Cells.Item(1, 9).Value2 = 7
$ws.Cells.Item(1, 10).Value2 = 8
$ws.Cells.Item(1, 11).Value2 = 9

# 3) Overwrite the data cells B2:K12 with the updated values
# Row 2
$ws.Cells.Item(2, 2).Value2 = 0.9419999999999999
$ws.Cells.Item(2, 3).Value2 = 0.873
$ws.Cells.Item(2, 4).Value2 = -2.993
$ws.Cells.Item(2, 5).Value2 = 1.018
$ws.Cells.Item(2, 6).Value2 = 7.254
$ws.Cells.Item(2, 7).Value2 = 0.141
$ws.Cells.Item(2, 8).Value2 = 0.127
$ws.Cells.Item(2, 9).Value2 = 0.577
$ws.Cells.Item(2, 10).Value2 = 412
$ws.Cells.Item(2, 11).Value2 = 1919
# Row 3
$ws.Cells.Item(3, 2).Value2 = 0.97
$ws.Cells.Item(3, 3).Value2 = 0.921
$ws.Cells.Item(3, 4).Value2 = -3.003
$ws.Cells.Item(3, 5).Value2 = 1.005
$ws.Cells.Item(3, 6).Value2 = 5.634
$ws.Cells.Item(3, 7).Value2 = 0.08500000000000001
$ws.Cells.Item(3, 8).Value2 = 0.079
$ws.Cells.Item(3, 9).Value2 = 0.549
$ws.Cells.Item(3, 10).Value2 = 378
$ws.Cells.Item(3, 11).Value2 = 1759
# Row 4
$ws.Cells.Item(4, 2).Value2 = 0.971
$ws.Cells.Item(4, 3).Value2 = 0.895
$ws.Cells.Item(4, 4).Value2 = -3.015
$ws.Cells.Item(4, 5).Value2 = 1.034
$ws.Cells.Item(4, 6).Value2 = 4.181
$ws.Cells.Item(4, 7).Value2 = 0.115
$ws.Cells.Item(4, 8).Value2 = 0.105
$ws.Cells.Item(4, 9).Value2 = 0.676
$ws.Cells.Item(4, 10).Value2 = 152
$ws.Cells.Item(4, 11).Value2 = 695
# Row 5
$ws.Cells.Item(5, 2).Value2 = 1.006
$ws.Cells.Item(5, 3).Value2 = 0.848
$ws.Cells.Item(5, 4).Value2 = -3.023
$ws.Cells.Item(5, 5).Value2 = 1.047
$ws.Cells.Item(5, 6).Value2 = 3.145
$ws.Cells.Item(5, 7).Value2 = 0.161
$ws.Cells.Item(5, 8).Value2 = 0.152
$ws.Cells.Item(5, 9).Value2 = 0.699
$ws.Cells.Item(5, 10).Value2 = 178
$ws.Cells.Item(5, 11).Value2 = 815
# Row 6
$ws.Cells.Item(6, 2).Value2 = 1.043
$ws.Cells.Item(6, 3).Value2 = 0.873
$ws.Cells.Item(6, 4).Value2 = -3.06
$ws.Cells.Item(6, 5).Value2 = 1.053
$ws.Cells.Item(6, 6).Value2 = 2.218
$ws.Cells.Item(6, 7).Value2 = 0.156
$ws.Cells.Item(6, 8).Value2 = 0.127
$ws.Cells.Item(6, 9).Value2 = 0.6830000000000001
$ws.Cells.Item(6, 10).Value2 = 206
$ws.Cells.Item(6, 11).Value2 = 942
# Row 7
$ws.Cells.Item(7, 2).Value2 = 1.038
$ws.Cells.Item(7, 3).Value2 = 0.918
$ws.Cells.Item(7, 4).Value2 = -3.019
$ws.Cells.Item(7, 5).Value2 = 1.022
$ws.Cells.Item(7, 6).Value2 = 1.335
$ws.Cells.Item(7, 7).Value2 = 0.095
$ws.Cells.Item(7, 8).Value2 = 0.082
$ws.Cells.Item(7, 9).Value2 = 0.645
$ws.Cells.Item(7, 10).Value2 = 209
$ws.Cells.Item(7, 11).Value2 = 957
# Row 8
$ws.Cells.Item(8, 2).Value2 = 1.026
$ws.Cells.Item(8, 3).Value2 = 0.946
$ws.Cells.Item(8, 4).Value2 = -2.976
$ws.Cells.Item(8, 5).Value2 = 0.999
$ws.Cells.Item(8, 6).Value2 = 0.576
$ws.Cells.Item(8, 7).Value2 = 0.065
$ws.Cells.Item(8, 8).Value2 = 0.054
$ws.Cells.Item(8, 9).Value2 = 0.648
$ws.Cells.Item(8, 10).Value2 = 178
$ws.Cells.Item(8, 11).Value2 = 812
# Row 9
$ws.Cells.Item(9, 2).Value2 = 1.004
$ws.Cells.Item(9, 3).Value2 = 0.981
$ws.Cells.Item(9, 4).Value2 = -2.976
$ws.Cells.Item(9, 5).Value2 = 0.993
$ws.Cells.Item(9, 6).Value2 = 0.091
$ws.Cells.Item(9, 7).Value2 = 0.031
$ws.Cells.Item(9, 8).Value2 = 0.024
$ws.Cells.Item(9, 9).Value2 = 0.636
$ws.Cells.Item(9, 10).Value2 = 175
$ws.Cells.Item(9, 11).Value2 = 796
# Row 10
$ws.Cells.Item(10, 2).Value2 = 1.002
$ws.Cells.Item(10, 3).Value2 = 0.978
$ws.Cells.Item(10, 4).Value2 = -2.976
$ws.Cells.Item(10, 5).Value2 = 0.994
$ws.Cells.Item(10, 6).Value2 = 0.08599999999999999
$ws.Cells.Item(10, 7).Value2 = 0.033
$ws.Cells.Item(10, 8).Value2 = 0.024
$ws.Cells.Item(10, 9).Value2 = 0.639
$ws.Cells.Item(10, 10).Value2 = 192
$ws.Cells.Item(10, 11).Value2 = 871
# Row 11
$ws.Cells.Item(11, 2).Value2 = 1.002
$ws.Cells.Item(11, 3).Value2 = 0.98
$ws.Cells.Item(11, 4).Value2 = -2.975
$ws.Cells.Item(11, 5).Value2 = 0.992
$ws.Cells.Item(11, 6).Value2 = 0.081
$ws.Cells.Item(11, 7).Value2 = 0.033
$ws.Cells.Item(11, 8).Value2 = 0.025
$ws.Cells.Item(11, 9).Value2 = 0.629
$ws.Cells.Item(11, 10).Value2 = 162
$ws.Cells.Item(11, 11).Value2 = 730
# Row 12
$ws.Cells.Item(12, 2).Value2 = 1.001
$ws.Cells.Item(12, 3).Value2 = 0.984
$ws.Cells.Item(12, 4).Value2 = -2.971
$ws.Cells.Item(12, 5).Value2 = 0.987
$ws.Cells.Item(12, 6).Value2 = 0.077
$ws.Cells.Item(12, 7).Value2 = 0.036
$ws.Cells.Item(12, 8).Value2 = 0.029
$ws.Cells.Item(12, 9).Value2 = 0.615
$ws.Cells.Item(12, 10).Value2 = 142
$ws.Cells.Item(12, 11).Value2 = 633

$excel.CutCopyMode = $false

Write-Output "edit applied"
